$d = $word.ActiveDocument

$ids = @("p011v_1", "p011v_2", "p011v_3", "p011v_4", "p011v_5")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false,
                             $true, 1, $false, $old, 2)
}
